$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.786.84"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "1.629.43"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  -0.91%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.92"
$ws.Range("E5").Value = "  -0.62%  "

$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("E7").Value = "  -0.93%  "

$ws.Range("E8").Value = "  -1.19%  "

$ws.Range("E9").Value = "  -0.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("E10").Value = "  +0.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.24"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Value = "1.856.81"
$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("D14").Value = "1.627.49"
$ws.Range("E14").Value = "  -1.44%  "

$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("D16").Value = "0.0₃0757"
$ws.Range("E16").Value = "  -1.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.75"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").Value = "25.796.47"
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.77"
$ws.Range("E21").Value = "  -1.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.91"
$ws.Range("E22").Value = "  -0.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.26"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.996"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.80"
$ws.Range("E25").Value = "  +0.95%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.18"
$ws.Range("E26").Value = "  +1.64%  "

$ws.Range("E27").Value = "  +1.58%  "

$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.50"
$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  -0.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0493"
$ws.Range("E31").Value = "  +0.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("E33").Value = "  -0.58%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.902"
$ws.Range("E36").Value = "  +0.54%  "

$ws.Range("D37").Value = "1.140.00"
$ws.Range("E37").Value = "  +1.98%  "

$ws.Range("E38").Value = "  -0.59%  "

$ws.Range("E39").Value = "  -0.99%  "

$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.994"
$ws.Range("E41").Value = "  -0.76%  "

$ws.Range("E42").Value = "  -1.48%  "

$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.31"
$ws.Range("E44").Value = "  +0.63%  "

$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").Value = "1.767.07"
$ws.Range("E46").Value = "  -0.53%  "

$ws.Range("E47").Value = "  +1.64%  "

$ws.Range("E48").Value = "  +0.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.47"
$ws.Range("E49").Value = "  +6.70%  "

$ws.Range("E50").Value = "  +1.97%  "

$ws.Range("E51").Value = "  -0.57%  "
